# "Updated Daq testing code"
# The pin-header labels in row 2 (top connector) and row 5 (bottom connector)
# were swapped for columns D, G, J, M, O (the SCL/SDA columns in T stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$swapCols = @("D", "G", "J", "M", "O")
foreach ($col in $swapCols) {
    $topCell = $ws.Range("$col" + "2")
    $bottomCell = $ws.Range("$col" + "5")
    $topValue = $topCell.Value()
    $bottomValue = $bottomCell.Value()
    $topCell.Value = $bottomValue
    $bottomCell.Value = $topValue
}

# Move / update the active selection to reflect where the editor left off.
$ws.Range("O11").Select() | Out-Null
